$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add the new sentence documenting Jeff's sound-file work, between the
#    existing "...and even more documentation." sentence and the following
#    "Outside of the listed duties..." sentence.
# ---------------------------------------------------------------------------
$quoteOpen  = [char]0x201C
$quoteClose = [char]0x201D

$oldDuties = " and even more documentation.  Outside of the listed duties, a smaller dut"
$newDuties = " and even more documentation." `
    + " Finally, Jeff implemented sound files to be played if the user selects the " `
    + $quoteOpen + "Help" + $quoteClose `
    + " button and also at the end of the game. These sound files are different for each scenario such as winning or losing." `
    + " Outside of the listed duties, a smaller dut"

$d.Content.Find.Execute($oldDuties, $true, $false, $false, $false, $false, $true, 1, $false, $newDuties, 2)

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark (Word's "last edit location" marker) so it
#    sits at the end of the paragraph we just edited -- right after
#    "...TA Oqi as needed." -- instead of its old spot in the "code-freeze"
#    paragraph further down.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$para = $d.Content
$para.Find.Execute("TA Oqi as needed.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endOfParaText = $para.End

# A zero-length range sitting exactly on the paragraph-end boundary can't be
# handed to Bookmarks.Add directly, so temporarily extend the paragraph with
# a placeholder character, anchor the bookmark just before it (now a safe,
# non-boundary position) and then remove the placeholder again.
$tail = $d.Range($endOfParaText, $endOfParaText)
$tail.InsertAfter("@")

$point = $d.Range($endOfParaText, $endOfParaText)
$d.Bookmarks.Add("_GoBack", $point)

$placeholder = $d.Range($endOfParaText, $endOfParaText + 1)
$placeholder.Text = ""
